$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 127.84615
$ws.Range("I55").Value = 122.5
$ws.Range("J55").Value = 136.4
$ws.Range("K55").Value = 122.5
$ws.Range("L55").Value = 136.4
$ws.Range("M55").Value = 91.5
$ws.Range("N55").Value = -564.4
# Row 64
$ws.Range("H64").Value = 4260.56
$ws.Range("I64").Value = 3999.6667
$ws.Range("J64").Value = 4407.3125
$ws.Range("K64").Value = 3999.6667
$ws.Range("L64").Value = 4407.3125
$ws.Range("M64").Value = -3751.6667
$ws.Range("N64").Value = -4903.3125
# Row 67
$ws.Range("H67").Value = 4260.56
$ws.Range("I67").Value = 3999.6667
$ws.Range("J67").Value = 4407.3125
$ws.Range("K67").Value = 3999.6667
$ws.Range("L67").Value = 4407.3125
$ws.Range("M67").Value = -3141.6667
$ws.Range("N67").Value = -6123.3125
# Row 76
$ws.Range("H76").Value = 5724.375
$ws.Range("I76").Value = 4998.5713
$ws.Range("K76").Value = 4998.5713
$ws.Range("M76").Value = -4683.5713
# Row 79
$ws.Range("H79").Value = 5724.375
$ws.Range("I79").Value = 4998.5713
$ws.Range("K79").Value = 4998.5713
$ws.Range("M79").Value = -3906.5713
# Row 137
$ws.Range("H137").Value = 412373.38
$ws.Range("I137").Value = 745876.4399999999
$ws.Range("J137").Value = 1908.0769
$ws.Range("K137").Value = 2237629.32
$ws.Range("L137").Value = 5724.2307
$ws.Range("M137").Value = -2235079.32
$ws.Range("N137").Value = -10824.2307

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22178.275
$ws.Range("I32").Value = 20153.467
$ws.Range("J32").Value = 34833.332
$ws.Range("K32").Value = 20153.467
$ws.Range("L32").Value = 34833.332
$ws.Range("M32").Value = -19866.467
$ws.Range("N32").Value = -35407.332
# Row 63
$ws.Range("H63").Value = 1192757
$ws.Range("I63").Value = 1853506.5
$ws.Range("J63").Value = 3408
$ws.Range("K63").Value = 1853506.5
$ws.Range("L63").Value = 3408
$ws.Range("M63").Value = -1852820.5
$ws.Range("N63").Value = -4780
# Row 66
$ws.Range("H66").Value = 1192757
$ws.Range("I66").Value = 1853506.5
$ws.Range("J66").Value = 3408
$ws.Range("K66").Value = 9267532.5
$ws.Range("L66").Value = 17040
$ws.Range("M66").Value = -9264100.5
$ws.Range("N66").Value = -23904
# Row 88
$ws.Range("H88").Value = 21733.092
$ws.Range("I88").Value = 1262.5
$ws.Range("J88").Value = 33430.57
$ws.Range("K88").Value = 1262.5
$ws.Range("L88").Value = 33430.57
$ws.Range("M88").Value = -856.5
$ws.Range("N88").Value = -34242.57
# Row 91
$ws.Range("H91").Value = 21733.092
$ws.Range("I91").Value = 1262.5
$ws.Range("J91").Value = 33430.57
$ws.Range("K91").Value = 1262.5
$ws.Range("L91").Value = 33430.57
$ws.Range("M91").Value = 141.5
$ws.Range("N91").Value = -36238.57

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 19984.55
$ws.Range("I31").Value = 18285.166
$ws.Range("J31").Value = 20712.857
$ws.Range("K31").Value = 18285.166
$ws.Range("L31").Value = 20712.857
$ws.Range("M31").Value = -17990.166
$ws.Range("N31").Value = -21302.857
# Row 34
$ws.Range("H34").Value = 19984.55
$ws.Range("I34").Value = 18285.166
$ws.Range("J34").Value = 20712.857
$ws.Range("K34").Value = 18285.166
$ws.Range("L34").Value = 20712.857
$ws.Range("M34").Value = -18083.166
$ws.Range("N34").Value = -21116.857
# Row 58
$ws.Range("H58").Value = 4235174.5
$ws.Range("I58").Value = 6255485.5
$ws.Range("J58").Value = 10887.546
$ws.Range("K58").Value = 6255485.5
$ws.Range("L58").Value = 10887.546
$ws.Range("M58").Value = -6255282.5
$ws.Range("N58").Value = -11293.546
# Row 62
$ws.Range("H62").Value = 4834.1665
$ws.Range("I62").Value = 4668.3335
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4668.3335
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4044.3335
$ws.Range("N62").Value = -6248
# Row 65
$ws.Range("H65").Value = 4834.1665
$ws.Range("I65").Value = 4668.3335
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 23341.6675
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -20221.6675
$ws.Range("N65").Value = -31240
# Row 132
$ws.Range("H132").Value = 10421751
$ws.Range("I132").Value = 17545316
$ws.Range("J132").Value = 10387.77
$ws.Range("K132").Value = 52635948
$ws.Range("L132").Value = 31163.31
$ws.Range("M132").Value = -52633418
$ws.Range("N132").Value = -36223.31
# Row 134
$ws.Range("H134").Value = 14040425
$ws.Range("I134").Value = 35716828
$ws.Range("J134").Value = 5610712
$ws.Range("K134").Value = 107150484
$ws.Range("L134").Value = 16832136
$ws.Range("M134").Value = -107147949
$ws.Range("N134").Value = -16837206
# Row 136
$ws.Range("H136").Value = 4235174.5
$ws.Range("I136").Value = 6255485.5
$ws.Range("J136").Value = 10887.546
$ws.Range("K136").Value = 18766456.5
$ws.Range("L136").Value = 32662.638
$ws.Range("M136").Value = -18763906.5
$ws.Range("N136").Value = -37762.638

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 946.7907
$ws.Range("I5").Value = 394
$ws.Range("K5").Value = 1182
$ws.Range("M5").Value = -1070
# Row 115
$ws.Range("H115").Value = 2704
$ws.Range("I115").Value = 642.6667
$ws.Range("J115").Value = 4250
$ws.Range("K115").Value = 1928.0001
$ws.Range("L115").Value = 12750
$ws.Range("M115").Value = -753.0001
$ws.Range("N115").Value = -15100
# Row 135
$ws.Range("H135").Value = 946.7907
$ws.Range("I135").Value = 394
$ws.Range("K135").Value = 3546
$ws.Range("M135").Value = -1011

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 30979.082
$ws.Range("I70").Value = 45512.832
$ws.Range("J70").Value = 4147.5386
$ws.Range("K70").Value = 45512.832
$ws.Range("L70").Value = 4147.5386
$ws.Range("M70").Value = -45242.832
$ws.Range("N70").Value = -4687.5386
# Row 73
$ws.Range("H73").Value = 30979.082
$ws.Range("I73").Value = 45512.832
$ws.Range("J73").Value = 4147.5386
$ws.Range("K73").Value = 45512.832
$ws.Range("L73").Value = 4147.5386
$ws.Range("M73").Value = -44576.832
$ws.Range("N73").Value = -6019.5386
# Row 80
$ws.Range("H80").Value = 2644.2222
$ws.Range("J80").Value = 2662.25
$ws.Range("L80").Value = 2662.25
$ws.Range("N80").Value = -4658.25
# Row 83
$ws.Range("H83").Value = 2644.2222
$ws.Range("J83").Value = 2662.25
$ws.Range("L83").Value = 13311.25
$ws.Range("N83").Value = -23295.25
# Row 141
$ws.Range("H141").Value = 44857.25
$ws.Range("J141").Value = 44857.25
$ws.Range("L141").Value = 44857.25
$ws.Range("N141").Value = -55217.25

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3852.7646
$ws.Range("I7").Value = 3785.5
$ws.Range("J7").Value = 4166.6665
$ws.Range("K7").Value = 3785.5
$ws.Range("L7").Value = 4166.6665
$ws.Range("M7").Value = -3673.5
$ws.Range("N7").Value = -4390.6665
# Row 126
$ws.Range("H126").Value = 3852.7646
$ws.Range("I126").Value = 3785.5
$ws.Range("J126").Value = 4166.6665
$ws.Range("K126").Value = 11356.5
$ws.Range("L126").Value = 12499.9995
$ws.Range("M126").Value = -8886.5
$ws.Range("N126").Value = -17439.9995
# Row 132
$ws.Range("H132").Value = 5027528.5
$ws.Range("I132").Value = 5745146.5
$ws.Range("J132").Value = 4200.6
$ws.Range("K132").Value = 17235439.5
$ws.Range("L132").Value = 12601.8
$ws.Range("M132").Value = -17232909.5
$ws.Range("N132").Value = -17661.8
# Row 136
$ws.Range("H136").Value = 4220.0215
$ws.Range("I136").Value = 7159.8096
$ws.Range("J136").Value = 1845.5769
$ws.Range("K136").Value = 21479.4288
$ws.Range("L136").Value = 5536.7307
$ws.Range("M136").Value = -18929.4288
$ws.Range("N136").Value = -10636.7307

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1257.1428
$ws.Range("I122").Value = 1257.1428
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3771.4284
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1321.4284
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 264188.88
$ws.Range("I132").Value = 23317.412
$ws.Range("J132").Value = 1271469.5
$ws.Range("K132").Value = 69952.236
$ws.Range("L132").Value = 3814408.5
$ws.Range("M132").Value = -67422.236
$ws.Range("N132").Value = -3819468.5
# Row 135
$ws.Range("H135").Value = 47501.766
$ws.Range("J135").Value = 45502.617
$ws.Range("L135").Value = 45502.617
$ws.Range("N135").Value = -55642.617
# Row 136
$ws.Range("H136").Value = 583134.1
$ws.Range("I136").Value = 1038160.94
$ws.Range("J136").Value = 1711
$ws.Range("K136").Value = 3114482.82
$ws.Range("L136").Value = 5133
$ws.Range("M136").Value = -3111932.82
$ws.Range("N136").Value = -10233
# Row 138
$ws.Range("H138").Value = 51666.668
$ws.Range("J138").Value = 51666.668
$ws.Range("L138").Value = 51666.668
$ws.Range("N138").Value = -61946.668
# Row 140
$ws.Range("H140").Value = 70015
$ws.Range("J140").Value = 70015
$ws.Range("L140").Value = 70015
$ws.Range("N140").Value = -80375
